# Metro 52 Zuid -> Noord: the run direction / timing for this trip was
# rechecked (commit: "check metro 52 and output all in QGIS"), so the whole
# A1:A32 schedule column is replaced with the corrected station/time list.
# (Time labels use the same "H:MM<nbsp>+<nbsp>N" convention, with U+00A0
# non-breaking spaces, as the rest of the workbook.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => (isTime, value). isTime cells keep the existing "h:mm" time style
# (s="1" in the OOXML); text cells are reset back to the default/Normal
# style so no leftover time-number-format sticks to them.
$rows = @(
    @{ Row = 1;  IsTime = $false; Value = "Station Zuid" },
    @{ Row = 2;  IsTime = $false; Value = "13:31 + 1" },
    @{ Row = 3;  IsTime = $false; Value = "13:38 + 1" },
    @{ Row = 4;  IsTime = $true;  Value = 0.57361111111111118 },
    @{ Row = 5;  IsTime = $false; Value = "Europaplein" },
    @{ Row = 6;  IsTime = $false; Value = "13:33 + 1" },
    @{ Row = 7;  IsTime = $false; Value = "13:41 + 1" },
    @{ Row = 8;  IsTime = $true;  Value = 0.57500000000000007 },
    @{ Row = 9;  IsTime = $false; Value = "De Pijp" },
    @{ Row = 10; IsTime = $false; Value = "13:35 + 1" },
    @{ Row = 11; IsTime = $false; Value = "13:43 + 1" },
    @{ Row = 12; IsTime = $true;  Value = 0.57638888888888895 },
    @{ Row = 13; IsTime = $false; Value = "Vijzelgracht" },
    @{ Row = 14; IsTime = $false; Value = "13:37 + 1" },
    @{ Row = 15; IsTime = $false; Value = "13:44 + 1" },
    @{ Row = 16; IsTime = $true;  Value = 0.57777777777777783 },
    @{ Row = 17; IsTime = $false; Value = "Rokin" },
    @{ Row = 18; IsTime = $false; Value = "13:39 + 1" },
    @{ Row = 19; IsTime = $false; Value = "13:46 + 1" },
    @{ Row = 20; IsTime = $true;  Value = 0.57916666666666672 },
    @{ Row = 21; IsTime = $false; Value = "Centraal Station" },
    @{ Row = 22; IsTime = $false; Value = "13:41 + 1" },
    @{ Row = 23; IsTime = $true;  Value = 0.57500000000000007 },
    @{ Row = 24; IsTime = $true;  Value = 0.58055555555555558 },
    @{ Row = 25; IsTime = $false; Value = "Noorderpark" },
    @{ Row = 26; IsTime = $false; Value = "13:43 + 1" },
    @{ Row = 27; IsTime = $true;  Value = 0.57638888888888895 },
    @{ Row = 28; IsTime = $true;  Value = 0.58194444444444449 },
    @{ Row = 29; IsTime = $false; Value = "Noord" },
    @{ Row = 30; IsTime = $false; Value = "13:45 + 1" },
    @{ Row = 31; IsTime = $true;  Value = 0.57847222222222217 },
    @{ Row = 32; IsTime = $true;  Value = 0.58333333333333337 }
)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r.Row, 1)
    if (-not $r.IsTime) {
        # Drop any inherited time-number-format so text cells come back
        # with the plain/default style (no explicit style index).
        $cell.Style = "Normal"
    }
    $cell.Value = $r.Value
}
